$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from NullFlavor")

# 1. Rename the include sheet.
$ws2.Name = "Include #0"

# 2. Bump the Version and Date metadata values.
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 3. Insert a new "Jurisdiction" metadata row right after "Contact" (row 10),
#    pushing Description/Purpose/Copyright/Immutable down by one row.
#    Values are shifted manually (instead of Rows.Insert) so the existing
#    cell style ("s=2") is reused rather than a new style being allocated.

# Capture the old values of rows 11-14 (Description, Purpose, Copyright, Immutable).
$oldRows = @()
for ($r = 11; $r -le 14; $r++) {
    $oldRows += ,@($ws1.Range("A$r").Value2, $ws1.Range("B$r").Value2)
}

# Row 15 is brand new - give it the same formatting as row 14 before writing into it.
$ws1.Range("A14:B14").Copy()
$ws1.Range("A15:B15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the captured values back, shifted down by one row (12..15).
for ($i = 0; $i -le 3; $i++) {
    $destRow = 12 + $i
    $ws1.Range("A$destRow").Value = $oldRows[$i][0]
    $bVal = $oldRows[$i][1]
    if ($bVal -eq $null -or $bVal -eq "") {
        $ws1.Range("B$destRow").ClearContents()
    } else {
        $ws1.Range("B$destRow").Value = $bVal
    }
}

# Finally, populate the freed-up row 11 with the new "Jurisdiction" property.
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").ClearContents()
